{"js": "// Map of old equation text -> new equation text, taken from the diff.\nconst replacements = [\n  [\"80\u00d728=2240\", \"14\u00d768=952\"],\n  [\"60\u00d761=3660\", \"73\u00d746=3358\"],\n  [\"77\u00d728=2156\", \"68\u00d759=4012\"],\n  [\"96\u00d726=2496\", \"25\u00d758=1450\"],\n  [\"64\u00d722=1408\", \"88\u00d791=8008\"],\n  [\"72\u00d799=7128\", \"27\u00d734=918\"],\n  [\"76\u00d763=4788\", \"63\u00d780=5040\"],\n  [\"78\u00d746=3588\", \"28\u00d732=896\"],\n  [\"72\u00d788=6336\", \"68\u00d717=1156\"],\n  [\"41\u00d776=3116\", \"41\u00d761=2501\"],\n  [\"38\u00d753=2014\", \"69\u00d745=3105\"],\n  [\"34\u00d713=442\", \"61\u00d777=4697\"],\n  [\"46\u00d756=2576\", \"47\u00d737=1739\"],\n  [\"55\u00d753=2915\", \"32\u00d711=352\"],\n  [\"20\u00d762=1240\", \"11\u00d748=528\"],\n  [\"65\u00d734=2210\", \"52\u00d747=2444\"],\n  [\"94\u00d797=9118\", \"92\u00d766=6072\"],\n  [\"73\u00d729=2117\", \"87\u00d757=4959\"],\n  [\"56\u00d749=2744\", \"65\u00d780=5200\"],\n  [\"69\u00d734=2346\", \"47\u00d799=4653\"],\n  [\"26\u00d725=650\", \"36\u00d727=972\"],\n  [\"29\u00d761=1769\", \"32\u00d762=1984\"],\n  [\"86\u00d726=2236\", \"47\u00d790=4230\"],\n  [\"54\u00d798=5292\", \"37\u00d777=2849\"],\n  [\"44\u00d794=4136\", \"52\u00d739=2028\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old equation text -> new equation text, taken from the diff.\n$replacements = @(\n    @(\"80\u00d728=2240\", \"14\u00d768=952\"),\n    @(\"60\u00d761=3660\", \"73\u00d746=3358\"),\n    @(\"77\u00d728=2156\", \"68\u00d759=4012\"),\n    @(\"96\u00d726=2496\", \"25\u00d758=1450\"),\n    @(\"64\u00d722=1408\", \"88\u00d791=8008\"),\n    @(\"72\u00d799=7128\", \"27\u00d734=918\"),\n    @(\"76\u00d763=4788\", \"63\u00d780=5040\"),\n    @(\"78\u00d746=3588\", \"28\u00d732=896\"),\n    @(\"72\u00d788=6336\", \"68\u00d717=1156\"),\n    @(\"41\u00d776=3116\", \"41\u00d761=2501\"),\n    @(\"38\u00d753=2014\", \"69\u00d745=3105\"),\n    @(\"34\u00d713=442\",  \"61\u00d777=4697\"),\n    @(\"46\u00d756=2576\", \"47\u00d737=1739\"),\n    @(\"55\u00d753=2915\", \"32\u00d711=352\"),\n    @(\"20\u00d762=1240\", \"11\u00d748=528\"),\n    @(\"65\u00d734=2210\", \"52\u00d747=2444\"),\n    @(\"94\u00d797=9118\", \"92\u00d766=6072\"),\n    @(\"73\u00d729=2117\", \"87\u00d757=4959\"),\n    @(\"56\u00d749=2744\", \"65\u00d780=5200\"),\n    @(\"69\u00d734=2346\", \"47\u00d799=4653\"),\n    @(\"26\u00d725=650\",  \"36\u00d727=972\"),\n    @(\"29\u00d761=1769\", \"32\u00d762=1984\"),\n    @(\"86\u00d726=2236\", \"47\u00d790=4230\"),\n    @(\"54\u00d798=5292\", \"37\u00d777=2849\"),\n    @(\"44\u00d794=4136\", \"52\u00d739=2028\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
